$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D21").Value = 0.54861111111111105
$ws.Range("D21").NumberFormat = "h:mm:ss"

$ws.Range("C22").Value = 0.57638888888888895
$ws.Range("C22").NumberFormat = "h:mm:ss"

$ws.Range("D22").Value = 0.76388888888888884
$ws.Range("D22").NumberFormat = "h:mm"

$ws.Range("C23").Value = 0.79166666666666663
$ws.Range("C23").NumberFormat = "h:mm:ss"

$ws.Range("D23").Select()
